# Add team record (Wins/Losses/Ties) columns to the COL_2021 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# copy header style from an existing header cell (AC1) so the new headers
# match the bold/centered/bordered look of the rest of row 1
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-47: same team record for every player row
$wins = 74
$losses = 87
$ties = 0

for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($row, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($row, 32).Value = $ties    # column AF = 32
}
